$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '64.966.02'
$ws.Range('E2').Value = '  +0.14%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.526.77'
$ws.Range('E3').Value = '  -0.61%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '592.82'
$ws.Range('E5').Value = '  -1.09%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '133.96'
$ws.Range('E6').Value = '  -1.95%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.525.21'
$ws.Range('E7').Value = '  -0.61%  '
$ws.Range('E8').Value = '  +0.18%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.491'
$ws.Range('E9').Value = '  -1.05%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.125'
$ws.Range('E10').Value = '  +1.42%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.12'
$ws.Range('E11').Value = '  +2.88%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.387'
$ws.Range('E12').Value = '  +0.13%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.120.81'
$ws.Range('E13').Value = '  -0.73%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.73'
$ws.Range('E14').Value = '  +2.14%  '
$ws.Range('E15').Value = '  -0.01%  '
$ws.Range('E16').Value = '  +0.50%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.519.77'
$ws.Range('E17').Value = '  -0.75%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '65.005.02'
$ws.Range('E18').Value = '  +0.33%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.18'
$ws.Range('E19').Value = '  +1.46%  '
$ws.Range('E20').Value = '  +0.69%  '
$ws.Range('E21').Value = '  -2.01%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '392.33'
$ws.Range('E22').Value = '  +1.01%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.581'
$ws.Range('E23').Value = '  +0.79%  '
$ws.Range('E24').Value = '  +0.91%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.667.46'
$ws.Range('E25').Value = '  -0.66%  '
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('E27').Value = '  -3.23%  '
$ws.Range('E28').Value = '  +0.85%  '
$ws.Range('E29').Value = '  +9.93%  '
$ws.Range('E30').Value = '  +0.07%  '
$ws.Range('E31').Value = '  -0.58%  '
$ws.Range('E32').Value = '  +0.22%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.527.71'
$ws.Range('E33').Value = '  -0.83%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '24.14'
$ws.Range('E34').Value = '  +0.77%  '
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('E36').Value = '  -0.25%  '
$ws.Range('E37').Value = '  +6.18%  '
$ws.Range('E38').Value = '  +3.35%  '
$ws.Range('E39').Value = '  +0.36%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '168.20'
$ws.Range('E40').Value = '  -1.13%  '
$ws.Range('E41').Value = '  +1.09%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.822'
$ws.Range('E42').Value = '  -0.47%  '
$ws.Range('E43').Value = '  +4.97%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '25.98'
$ws.Range('E44').Value = '  -3.34%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '42.93'
$ws.Range('E45').Value = '  +0.67%  '
$ws.Range('E46').Value = '  +0.04%  '
$ws.Range('E47').Value = '  -0.04%  '
$ws.Range('E48').Value = '  +1.18%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.91'
$ws.Range('E49').Value = '  +0.02%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.414.25'
$ws.Range('E50').Value = '  -1.40%  '
$ws.Range('E51').Value = '  +5.54%  '
